# Applies three changes described in the commit:
#   1. Split the "Hinzufügen Button unter Combox..." bullet into three runs
#      ("Hinzufügen Button " / "unter" / " Combox und...") and move the
#      "_GoBack" bookmark so that it sits right after "unter".
#   2. Strike through the "Gesamt in Diagram als erste Wahl" bullet.
#   3. Remove the "_GoBack" bookmark that used to sit at the end of the
#      "Who you are Mihail" paragraph (it has moved, see step 1).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 & 3: relocate the "_GoBack" bookmark into the "Hinzufügen Button"
# bullet, splitting that paragraph's single run into three runs in the
# process.
# ---------------------------------------------------------------------------

# Locate the target paragraph without depending on a hard-coded paragraph
# index.
$target = $d.Content
$target.Find.Execute("Hinzufügen Button unter Combox und") | Out-Null
$paraStart = $target.Start

# Remove the bookmark from its old location (end of "Who you are Mihail").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Create a throw-away bookmark just to force Word to split the run after
# "Hinzufügen Button " (18 characters in, trailing space included).
$split1 = $paraStart + 18
$d.Bookmarks.Add("TempSplit1", $d.Range($split1, $split1))

# Re-create "_GoBack" right after "unter" (23 characters in), which both
# splits off the "unter" run and plants the bookmark in its new home.
$split2 = $paraStart + 23
$d.Bookmarks.Add("_GoBack", $d.Range($split2, $split2))

# Drop the helper bookmark now that it has done its job.
$d.Bookmarks("TempSplit1").Delete()

# ---------------------------------------------------------------------------
# Change 2: strike through "Gesamt in Diagram als erste Wahl".
# ---------------------------------------------------------------------------

$target2 = $d.Content
$target2.Find.Execute("Gesamt in Diagram als erste Wahl") | Out-Null
# Expand to the whole paragraph (incl. paragraph mark) so the strike
# formatting also lands on the paragraph mark's run properties.
$para2 = $target2.Paragraphs(1).Range
$para2.Font.StrikeThrough = $true
